$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values from diff.
# Numeric-looking text values (prices) use a leading apostrophe via Formula
# so Excel stores them as text (quotePrefix), matching the original inlineStr cells
# and preserving exact formatting (trailing zeros, precision).
$ws.Range("D2").Formula = "'244.30"
$ws.Range("D3").Formula = "'21.82"
$ws.Range("D5").Formula = "'0.05996"
$ws.Range("D6").Formula = "'3.391"
$ws.Range("D7").Formula = "'0.8166"
$ws.Range("D8").Formula = "'0.9557"
$ws.Range("D9").Formula = "'0.1425"
$ws.Range("D10").Formula = "'0.07436"
$ws.Range("D11").Formula = "'0.03305"
$ws.Range("D12").Formula = "'0.03055"
$ws.Range("D13").Formula = "'0.09407"
$ws.Range("D14").Formula = "'4.004"
$ws.Range("D15").Formula = "'0.001600"
$ws.Range("D16").Formula = "'0.04814"
$ws.Range("D18").Formula = "'0.005505"
$ws.Range("D19").Formula = "'0.004150"
$ws.Range("D20").Formula = "'0.0009893"
$ws.Range("B21").Value = 'LEO'
$ws.Range("C21").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D21").Formula = "'3.673"
$ws.Range("E21").Value = '20LEOLEO'
$ws.Range("B22").Value = 'KuCoinToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D22").Formula = "'6.425"
$ws.Range("E22").Value = '21KuCoinTokenKCS'
$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").Formula = "'2.189"
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("B24").Value = 'BitpandaEcosystemToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D24").Formula = "'0.3253"
$ws.Range("E24").Value = '23BitpandaEcosystemTokenBEST'
$ws.Range("B25").Value = 'ProBitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D25").Formula = "'0.1331"
$ws.Range("E25").Value = '24ProBitTokenPROB'
$ws.Range("B26").Value = 'NitroEx'
$ws.Range("C26").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D26").Formula = "'0.00007002"
$ws.Range("E26").Value = '25NitroExNTX'
$ws.Range("D40").Formula = "'0.04000"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Formula = "'0.1073"
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Formula = "'0.002721"
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Formula = "'0.003047"
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("D44").Formula = "'0.005799"
$ws.Range("D45").Formula = "'0.00005126"
$ws.Range("D47").Formula = "'0.8603"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINBestin24h'
$ws.Range("D48").Formula = "'0.005111"
